# BoardTestLayout.xlsx edit: "Failing adjacency tests complete"
#
# Recolors a handful of board-layout test cells (swapping which fill/style
# they use) and relabels several "W" (walkway) cells with more specific
# walkway-scenario codes (W1, W2, W4, W6, WIR, Wir) that were previously all
# just generic "W" placeholders. Also updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122  # xlPasteFormats

# ---------------------------------------------------------------------
# 1) Style-only changes: reuse an existing cell's exact format (fill) by
#    copy/paste-special so we don't create redundant duplicate styles.
# ---------------------------------------------------------------------

# target style 8 (green fill) <- donor L7
[void]$ws.Cells.Item(7,12).Copy()
[void]$ws.Cells.Item(1,8).PasteSpecial($xlPasteFormats)     # H1: style 9 -> 8

# target style 2 (gray/theme fill) <- donor L1
[void]$ws.Cells.Item(1,12).Copy()
[void]$ws.Cells.Item(2,11).PasteSpecial($xlPasteFormats)    # K2: style 8 -> 2
[void]$ws.Cells.Item(5,9).PasteSpecial($xlPasteFormats)     # I5: style 9 -> 2
[void]$ws.Cells.Item(8,12).PasteSpecial($xlPasteFormats)    # L8: style 9 -> 2
[void]$ws.Cells.Item(12,6).PasteSpecial($xlPasteFormats)    # F12: style 9 -> 2

# target style 3 (light fill) <- donor B1
[void]$ws.Cells.Item(1,2).Copy()
[void]$ws.Cells.Item(3,4).PasteSpecial($xlPasteFormats)     # D3: style 9 -> 3

# target style 4 (light blue fill) <- donor M2
[void]$ws.Cells.Item(2,13).Copy()
[void]$ws.Cells.Item(11,16).PasteSpecial($xlPasteFormats)   # P11: style 5 -> 4
[void]$ws.Cells.Item(20,18).PasteSpecial($xlPasteFormats)   # R20: style 9 -> 4

# target style 5 (blue fill) <- donor M3
[void]$ws.Cells.Item(3,13).Copy()
[void]$ws.Cells.Item(13,18).PasteSpecial($xlPasteFormats)   # R13: style 4 -> 5
[void]$ws.Cells.Item(17,4).PasteSpecial($xlPasteFormats)    # D17: style 6 -> 5

# target style 9 (red fill) <- donor D7  (K1 also gets a new value below)
[void]$ws.Cells.Item(7,4).Copy()
[void]$ws.Cells.Item(1,11).PasteSpecial($xlPasteFormats)    # K1: style 2 -> 9

# ---------------------------------------------------------------------
# 2) New style: plain white (theme "Background 1", no tint) fill, applied
#    to M19. Excel doesn't already have this exact fill defined, so this
#    allocates one new fill + cellXfs entry.
# ---------------------------------------------------------------------
$m19 = $ws.Cells.Item(19,13)
$m19.Interior.Pattern = 1                 # xlSolid
$m19.Interior.ThemeColor = 2              # xlThemeColorLight1 -> theme "0"
$m19.Interior.TintAndShade = 0

# ---------------------------------------------------------------------
# 3) Relabel walkway test cells. New shared-string values are created in
#    this order so they land at the same indices as the authoritative
#    edit: W1, W6, W4, W2, WIR, Wir.
# ---------------------------------------------------------------------
$ws.Cells.Item(11,13).Value = "W1"    # M11
$ws.Cells.Item(8,21).Value  = "W6"    # U8
$ws.Cells.Item(1,11).Value  = "W4"    # K1
$ws.Cells.Item(1,18).Value  = "W2"    # R1
$ws.Cells.Item(18,5).Value  = "WIR"   # E18
$ws.Cells.Item(11,15).Value = "Wir"   # O11
$ws.Cells.Item(22,6).Value  = "W1"    # F22 (reuse)
$ws.Cells.Item(15,1).Value  = "W4"    # A15 (reuse)
$ws.Cells.Item(10,5).Value  = "W2"    # E10 (reuse)

# ---------------------------------------------------------------------
# 4) Move the active selection to D3 (was I7).
# ---------------------------------------------------------------------
[void]$ws.Range("D3").Select()
